# Update cryptocurrency price/volume data (Price in column D, Volume(1h) in column E)
# for rows 2-51 of Sheet1, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.027.01'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '3.385.17'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.74'
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.13'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.382.73'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.62'
$ws.Range("E10").Value = '  +1.55%  '
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").Value = '3.961.00'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.23'
$ws.Range("E15").Value = '  +3.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  -2.85%  '
$ws.Range("D17").Value = '3.382.28'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '61.140.39'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.82'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.46'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '377.44'
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("D24").Value = '3.526.07'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.27'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("E28").Value = '  +11.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.173'
$ws.Range("E29").Value = '  +8.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.52'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.72'
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("E36").Value = '  -3.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.56'
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.88'
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.57'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0761'
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -3.20%  '
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.64'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.07'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("D48").Value = '2.481.43'
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.26'
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.80'
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.40'
$ws.Range("E51").Value = '  +3.66%  '
